$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: push the existing 2022-Q3 summary row down to row 3,
#    and write the new 2022-Q4 summary into row 2.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Cells.Item(3,1).Value = 1

$wsTotal.Cells.Item(2,2).Value = "2022-Q4"
$wsTotal.Cells.Item(2,3).Value = 5
$wsTotal.Cells.Item(2,4).Value = 0.17

# ---------------------------------------------------------------------
# 2) The old "2022-Q3" fund-holdings sheet becomes the new, separate
#    "2022-Q3" tab (placed right after it) so its original data is kept
#    untouched, while the original tab is repurposed to hold the new
#    2022-Q4 fund-holdings table.
# ---------------------------------------------------------------------
$wsQ3old = $wb.Worksheets.Item(2)

$wsQ4 = $wsQ3old
$wsQ4.Name = "2022-Q4"

$wsQ3new = $wb.Worksheets.Add($null, $wsQ4)
$wsQ3new.Name = "2022-Q3"
$wsQ4.Range("A1:H2").Copy($wsQ3new.Range("A1"))
$wsQ3new.Cells.Item(1,1).Clear()

# Clear out the old single-fund row so the new table can be written.
$wsQ4.Cells.Clear()

# Header row (bold/bordered style, matching the workbook's existing
# header style used on row 1 of the "总计" sheet).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $wsTotal.Range("B1").Copy($wsQ4.Cells.Item(1, $col))
    $wsQ4.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$rows = @(
    @("004671", "中融核心成长灵活配置混合",     "1.14", "65.86", "4.37", "0.0498", 4),
    @("010009", "中融成长优选混合C",             "1.05", "60.81", "3.79", "0.0398", 7),
    @("008422", "中融研发创新混合A",             "1.07", "61.65", "3.58", "0.0383", 5),
    @("008423", "中融研发创新混合C",             "0.62", "61.65", "3.58", "0.0222", 5),
    @("010008", "中融成长优选混合A",             "0.57", "60.81", "3.79", "0.0216", 7)
)

$r = 2
$idx = 0
foreach ($row in $rows) {
    $wsTotal.Range("A2").Copy($wsQ4.Cells.Item($r, 1))
    $wsQ4.Cells.Item($r, 1).Value = $idx

    $wsQ4.Cells.Item($r, 2).Value = "'" + $row[0]
    $wsQ4.Cells.Item($r, 2).Style = "Normal"
    $wsQ4.Cells.Item($r, 3).Value = $row[1]
    $wsQ4.Cells.Item($r, 4).Value = "'" + $row[2]
    $wsQ4.Cells.Item($r, 4).Style = "Normal"
    $wsQ4.Cells.Item($r, 5).Value = "'" + $row[3]
    $wsQ4.Cells.Item($r, 5).Style = "Normal"
    $wsQ4.Cells.Item($r, 6).Value = "'" + $row[4]
    $wsQ4.Cells.Item($r, 6).Style = "Normal"
    $wsQ4.Cells.Item($r, 7).Value = "'" + $row[5]
    $wsQ4.Cells.Item($r, 7).Style = "Normal"
    $wsQ4.Cells.Item($r, 8).Value = $row[6]

    $r++
    $idx++
}

# Restore the originally active tab ("总计").
$wsTotal.Activate()
